$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2784.762
$ws.Range("J17").Value = 2784.762
$ws.Range("L17").Value = 8354.286
$ws.Range("N17").Value = -8690.286

$ws.Range("H40").Value = 3622.5
$ws.Range("I40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("M40").Value = -1325

$ws.Range("H42").Value = 2892.5715
$ws.Range("J42").Value = 5023.25
$ws.Range("L42").Value = 15069.75
$ws.Range("N42").Value = -15529.75

$ws.Range("H69").Value = 2000
$ws.Range("J69").Value = 2000
$ws.Range("L69").Value = 6000
$ws.Range("N69").Value = -7748

$ws.Range("H72").Value = 2000
$ws.Range("J72").Value = 2000
$ws.Range("L72").Value = 18000
$ws.Range("N72").Value = -26736

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

$ws.Range("H132").Value = 2374.1304
$ws.Range("I132").Value = 2254.7727
$ws.Range("K132").Value = 6764.3181
$ws.Range("M132").Value = -4234.3181

$ws.Range("H137").Value = 1722.1
$ws.Range("J137").Value = 2146
$ws.Range("L137").Value = 6438
$ws.Range("N137").Value = -11538

$ws.Range("H138").Value = 3599.4194
$ws.Range("I138").Value = 2982.3333
$ws.Range("K138").Value = 8946.999899999999
$ws.Range("M138").Value = -3806.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3800
$ws.Range("I2").Value = 3800
$ws.Range("K2").Value = 3800
$ws.Range("M2").Value = -3687

$ws.Range("H45").Value = 2219.8
$ws.Range("I45").Value = 1100
$ws.Range("J45").Value = 2499.75
$ws.Range("K45").Value = 1100
$ws.Range("L45").Value = 2499.75
$ws.Range("M45").Value = -723
$ws.Range("N45").Value = -3253.75

$ws.Range("H63").Value = 3530

$ws.Range("H66").Value = 3530

$ws.Range("H116").Value = 3800
$ws.Range("I116").Value = 3800
$ws.Range("K116").Value = 3800
$ws.Range("M116").Value = -1506

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3800
$ws.Range("I3").Value = 3800
$ws.Range("K3").Value = 3800
$ws.Range("M3").Value = -3686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 650
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H16").Value = 1488.1666
$ws.Range("I16").Value = 1488.1666
$ws.Range("K16").Value = 1488.1666
$ws.Range("M16").Value = -1201.1666

$ws.Range("H86").Value = 9136.091
$ws.Range("I86").Value = 8289.888999999999
$ws.Range("K86").Value = 8289.888999999999
$ws.Range("M86").Value = -7166.888999999999

$ws.Range("H89").Value = 9136.091
$ws.Range("I89").Value = 8289.888999999999
$ws.Range("K89").Value = 41449.44499999999
$ws.Range("M89").Value = -35833.44499999999

$ws.Range("H113").Value = 1488.1666
$ws.Range("I113").Value = 1488.1666
$ws.Range("K113").Value = 1488.1666
$ws.Range("M113").Value = 681.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1945.1904
$ws.Range("I4").Value = 1796.8125
$ws.Range("J4").Value = 2420
$ws.Range("K4").Value = 5390.4375
$ws.Range("L4").Value = 7260
$ws.Range("M4").Value = -5278.4375
$ws.Range("N4").Value = -7484

$ws.Range("H7").Value = 313.1
$ws.Range("I7").Value = 190
$ws.Range("J7").Value = 395.16666
$ws.Range("K7").Value = 570
$ws.Range("L7").Value = 1185.49998
$ws.Range("M7").Value = -458
$ws.Range("N7").Value = -1409.49998

$ws.Range("H139").Value = 4475
$ws.Range("I139").Value = 600
$ws.Range("K139").Value = 1800
$ws.Range("M139").Value = 3340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = 29

$ws.Range("H25").Value = 22000
$ws.Range("I25").Value = 4000
$ws.Range("J25").Value = 40000
$ws.Range("K25").Value = 4000
$ws.Range("L25").Value = 40000
$ws.Range("M25").Value = -3471
$ws.Range("N25").Value = -41058

$ws.Range("H113").Value = 2508.4
$ws.Range("I113").Value = 2514.3333
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 2514.3333
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = -344.3332999999998
$ws.Range("N113").Value = -6839.5

$ws.Range("H122").Value = 2332.8823
$ws.Range("I122").Value = 1440.1428
$ws.Range("K122").Value = 4320.428400000001
$ws.Range("M122").Value = -1870.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1235.7858
$ws.Range("I22").Value = 1229.9
$ws.Range("J22").Value = 1250.5
$ws.Range("K22").Value = 1229.9
$ws.Range("L22").Value = 1250.5
$ws.Range("M22").Value = -934.9000000000001
$ws.Range("N22").Value = -1840.5

$ws.Range("H27").Value = 1235.7858
$ws.Range("I27").Value = 1229.9
$ws.Range("J27").Value = 1250.5
$ws.Range("K27").Value = 1229.9
$ws.Range("L27").Value = 1250.5
$ws.Range("M27").Value = -1122.9
$ws.Range("N27").Value = -1464.5

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H63").Value = 29600
$ws.Range("J63").Value = 29600
$ws.Range("L63").Value = 29600
$ws.Range("N63").Value = -31098

$ws.Range("H66").Value = 29600
$ws.Range("J66").Value = 29600
$ws.Range("L66").Value = 88800
$ws.Range("N66").Value = -96288

$ws.Range("H68").Value = 4109.4
$ws.Range("I68").Value = 4013
$ws.Range("K68").Value = 4013
$ws.Range("M68").Value = -3264

$ws.Range("H71").Value = 4109.4
$ws.Range("I71").Value = 4013
$ws.Range("K71").Value = 20065
$ws.Range("M71").Value = -16321

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5097.364
$ws.Range("I62").Value = 4045.25
$ws.Range("J62").Value = 5698.5713
$ws.Range("K62").Value = 4045.25
$ws.Range("L62").Value = 5698.5713
$ws.Range("M62").Value = -3421.25
$ws.Range("N62").Value = -6946.5713

$ws.Range("H65").Value = 5097.364
$ws.Range("I65").Value = 4045.25
$ws.Range("J65").Value = 5698.5713
$ws.Range("K65").Value = 20226.25
$ws.Range("L65").Value = 28492.8565
$ws.Range("M65").Value = -17106.25
$ws.Range("N65").Value = -34732.85649999999

Write-Output "Applied all Marilith_Profits updates"
